$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 - use same formatting as other header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I and J, rows 2-20
$data = @{
    2  = @(3, 6)
    3  = @(6, 7)
    4  = @(2, 4)
    5  = @(8, 8)
    6  = @(3, 5)
    7  = @(8, 8)
    8  = @(7, 8)
    9  = @(5, 6)
    10 = @(6, 7)
    11 = @(6, 8)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(7, 8)
    16 = @(1, 2)
    17 = @(5, 5)
    18 = @(4, 5)
    19 = @(4, 4)
    20 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
